$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.189.52'
$ws.Cells.Item(2, 5).Value = '  +1.33%  '
$ws.Cells.Item(3, 4).Value = '3.138.69'
$ws.Cells.Item(3, 5).Value = '  +3.47%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '579.82'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.14%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '174.41'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +3.88%  '
$ws.Cells.Item(8, 4).Value = '3.134.92'
$ws.Cells.Item(8, 5).Value = '  +3.45%  '
$ws.Cells.Item(9, 5).Value = '  +0.53%  '
$ws.Cells.Item(10, 5).Value = '  -2.22%  '
$ws.Cells.Item(11, 5).Value = '  +2.48%  '
$ws.Cells.Item(12, 5).Value = '  -0.81%  '
$ws.Cells.Item(13, 5).Value = '  +1.07%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '37.43'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.17%  '
$ws.Cells.Item(15, 5).Value = '  -0.97%  '
$ws.Cells.Item(16, 4).Value = '3.663.75'
$ws.Cells.Item(16, 5).Value = '  +3.57%  '
$ws.Cells.Item(17, 4).Value = '67.180.92'
$ws.Cells.Item(17, 5).Value = '  +1.30%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.19'
$cell.Style = "Normal"
$ws.Cells.Item(19, 4).Value = '3.141.95'
$ws.Cells.Item(19, 5).Value = '  +3.60%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.15'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.41%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '489.43'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +4.97%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.717'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.78%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.71'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +4.19%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '84.29'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.48%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.25'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +4.21%  '
$ws.Cells.Item(26, 5).Value = '  +2.95%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.03'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.20%  '
$ws.Cells.Item(28, 5).Value = '  +0.07%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.97'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -3.00%  '
$ws.Cells.Item(30, 5).Value = '  -0.98%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.69'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.28%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '29.07'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +3.13%  '
$ws.Cells.Item(33, 4).Value = '0.0₃0999'
$ws.Cells.Item(33, 5).Value = '  +0.78%  '
$ws.Cells.Item(34, 5).Value = '  -2.90%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.01%  '
$ws.Cells.Item(36, 5).Value = '  +1.50%  '
$ws.Cells.Item(37, 5).Value = '  -0.13%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '47.45'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.99%  '
$ws.Cells.Item(39, 5).Value = '  +2.63%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '50.07'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +1.10%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.312'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -1.14%  '
$ws.Cells.Item(42, 5).Value = '  +1.52%  '
$ws.Cells.Item(43, 5).Value = '  +0.20%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.78'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -2.18%  '
$ws.Cells.Item(45, 4).Value = '2.848.79'
$ws.Cells.Item(45, 5).Value = '  +4.88%  '
$ws.Cells.Item(46, 5).Value = '  -0.41%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '384.31'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.08%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '135.52'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +1.34%  '
$ws.Cells.Item(49, 5).Value = '  +0.00%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '24.90'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +1.70%  '
$ws.Cells.Item(51, 5).Value = '  -0.35%  '
